# Generate Report for Handoff
# Replaces the localized file's generated GUID-based name throughout the
# workbook and refreshes the handoff/handback timestamps, clearing out the
# stale "Latest Target File" / "Latest Handback File" + datetime values
# that no longer apply now that a fresh handoff cycle has started.

$wb = $excel.ActiveWorkbook

$oldName = "ebe840a9-4e14-4a32-93d4-fe2121ecce10"
$newName = "31c0d168-3e4c-4ef5-9b0d-c65dc5a96396"

$oldMd = "$oldName.md"
$newMd = "$newName.md"
$oldMdPath = "e2e\$oldName.md"
$newMdPath = "e2e\$newName.md"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value2 = $newMd
$wsOverview.Range("B2").Value2 = $newMdPath
$wsOverview.Range("G2").Value2 = "2016-08-18 15:05:38"

foreach ($h in @($wsOverview.Hyperlinks)) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$2') {
        $h.TextToDisplay = $newMdPath
    }
}

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value2 = $newMd
$wsZhCn.Range("G2").Value2 = "$newName.a7a69d61d583339ddd8b0d706cb0f064f8e0ca28.zh-cn.xlf"
$wsZhCn.Range("H2").Value2 = "2016-08-18 15:05:33"
$wsZhCn.Range("K2").Value2 = "0001-01-01 00:00:00"

# "Latest Target File" (I2) and "Latest Handback File" (J2) no longer
# apply to this handoff cycle, so they're cleared out.
foreach ($h in @($wsZhCn.Hyperlinks)) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $newMd
    } elseif ($addr -eq '$I$2') {
        $h.Delete()
    }
}
$wsZhCn.Range("I2").Value2 = ""
$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("J2").Value2 = ""

$wsZhCn.Columns.Item(9).ColumnWidth = 17.83
$wsZhCn.Columns.Item(10).ColumnWidth = 20.83

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value2 = $newMd
$wsDeDe.Range("G2").Value2 = "$newName.a7a69d61d583339ddd8b0d706cb0f064f8e0ca28.de-de.xlf"
$wsDeDe.Range("H2").Value2 = "2016-08-18 15:05:38"
$wsDeDe.Range("K2").Value2 = "0001-01-01 00:00:00"

foreach ($h in @($wsDeDe.Hyperlinks)) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') {
        $h.TextToDisplay = $newMd
    } elseif ($addr -eq '$I$2') {
        $h.Delete()
    }
}
$wsDeDe.Range("I2").Value2 = ""
$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("J2").Value2 = ""

$wsDeDe.Columns.Item(9).ColumnWidth = 17.83
$wsDeDe.Columns.Item(10).ColumnWidth = 20.83
